$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D, E changed
$ws.Range('D2').Value = '70.031.24'
$ws.Range('E2').Value = '  -0.47%  '

# Row 3: D, E changed
$ws.Range('D3').Value = '3.741.41'
$ws.Range('E3').Value = '  -0.90%  '

# Row 4: D, E changed
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.19%  '

# Row 5: D, E changed
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '619.11'
$ws.Range('E5').Value = '  +0.08%  '

# Row 6: D, E changed
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.52'
$ws.Range('E6').Value = '  +2.26%  '

# Row 7: D, E changed
$ws.Range('D7').Value = '3.746.78'
$ws.Range('E7').Value = '  -0.78%  '

# Row 8: D, E changed
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.17%  '

# Row 9: D, E changed
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.532'
$ws.Range('E9').Value = '  -2.87%  '

# Row 10: D, E changed
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.166'
$ws.Range('E10').Value = '  -2.45%  '

# Row 11: D, E changed
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.34'
$ws.Range('E11').Value = '  -1.16%  '

# Row 12: D, E changed
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.482'
$ws.Range('E12').Value = '  -4.55%  '

# Row 13: D, E changed
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.98'
$ws.Range('E13').Value = '  -1.28%  '

# Row 14: D, E changed
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000254'
$ws.Range('E14').Value = '  -2.16%  '

# Row 15: D, E changed
$ws.Range('D15').Value = '4.353.63'
$ws.Range('E15').Value = '  -1.26%  '

# Row 16: D, E changed
$ws.Range('D16').Value = '3.734.10'
$ws.Range('E16').Value = '  -1.54%  '

# Row 17: D, E changed
$ws.Range('D17').Value = '70.040.09'
$ws.Range('E17').Value = '  -0.75%  '

# Row 18: E changed
$ws.Range('E18').Value = '  -2.07%  '

# Row 19: D, E changed
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.58'
$ws.Range('E19').Value = '  -0.38%  '

# Row 20: D, E changed
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '503.32'
$ws.Range('E20').Value = '  -4.16%  '

# Row 21: D, E changed
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.37'
$ws.Range('E21').Value = '  -3.44%  '

# Row 22: D, E changed
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.31'
$ws.Range('E22').Value = '  -1.11%  '

# Row 23: D, E changed
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.722'
$ws.Range('E23').Value = '  -3.20%  '

# Row 24: D, E changed
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.55'
$ws.Range('E24').Value = '  +2.73%  '

# Row 25: D, E changed
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.64'
$ws.Range('E25').Value = '  -1.65%  '

# Row 26: D, E changed
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.98'
$ws.Range('E26').Value = '  -4.19%  '

# Row 27: D, E changed
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.20'
$ws.Range('E27').Value = '  +2.24%  '

# Row 28: D, E changed
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000132'
$ws.Range('E28').Value = '  +7.13%  '

# Row 29: E changed
$ws.Range('E29').Value = '  +0.28%  '

# Row 30: D, E changed
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.47'
$ws.Range('E30').Value = '  -1.69%  '

# Row 31: D, E changed
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.92'
$ws.Range('E31').Value = '  +0.51%  '

# Row 32: D, E changed
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.90'
$ws.Range('E32').Value = '  -0.89%  '

# Row 33: D, E changed
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '30.46'
$ws.Range('E33').Value = '  -5.38%  '

# Row 34: D, E changed
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.115'
$ws.Range('E34').Value = '  -0.44%  '

# Row 35: D, E changed
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.26%  '

# Row 36: E changed
$ws.Range('E36').Value = '  +0.23%  '

# Row 37: D, E changed
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.09'
$ws.Range('E37').Value = '  -1.35%  '

# Row 38: B, C, D, E changed
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.347'
$ws.Range('E38').Value = '  +1.16%  '

# Row 39: B, C, D, E changed
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.139'
$ws.Range('E39').Value = '  +3.87%  '

# Row 40: D, E changed
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.14'
$ws.Range('E40').Value = '  +12.98%  '

# Row 41: D, E changed
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.08'
$ws.Range('E41').Value = '  -5.17%  '

# Row 42: D, E changed
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '50.04'
$ws.Range('E42').Value = '  -2.71%  '

# Row 43: D, E changed
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '427.18'
$ws.Range('E43').Value = '  -0.50%  '

# Row 44: D, E changed
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '44.30'
$ws.Range('E44').Value = '  -0.68%  '

# Row 45: D, E changed
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.60'
$ws.Range('E45').Value = '  -3.23%  '

# Row 46: D, E changed
$ws.Range('D46').Value = '2.956.90'
$ws.Range('E46').Value = '  -6.20%  '

# Row 47: D, E changed
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0361'
$ws.Range('E47').Value = '  -2.09%  '

# Row 48: B, C, D, E changed
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '27.18'
$ws.Range('E48').Value = '  -2.30%  '

# Row 49: B, C, D, E changed
$ws.Range('B49').Value = 'USDe'
$ws.Range('C49').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -0.08%  '

# Row 50: E changed
$ws.Range('E50').Value = '  -2.66%  '

# Row 51: D, E changed
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.47'
$ws.Range('E51').Value = '  -2.57%  '
